# convocazioni stagione 126 - 4
# Swap columns A and B (TagTelegram / Nome) for rows 52..143, then update
# the sheet view (scroll position / zoom / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 52; $r -le 143; $r++) {
    $colA = $ws.Cells.Item($r, 1)
    $colB = $ws.Cells.Item($r, 2)
    $a = $colA.Value()
    $b = $colB.Value()

    if ($b -eq $null) {
        # Row had only column A populated (B empty) -> after swap B holds
        # the value and A becomes empty (cell removed entirely).
        $colB.Value = $a
        $colA.Clear()
    } else {
        $colA.Value = $b
        $colB.Value = $a
    }
}

# Update the view: scrolled position, zoom level and active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 134
$win.ScrollColumn = 1
$win.Zoom = 143
$ws.Range("F142").Select() | Out-Null
